$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = " -"
$ws.Range("L2").Value = "-"
$ws.Range("N2").Value = "-"

$ws.Range("L3").Select()
